$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2071823204419889
$ws.Range("C2").Value = 0.5359116022099447
$ws.Range("J2").Value = 0.01933701657458563
$ws.Range("P2").Value = 0.1574585635359116
$ws.Range("S2").Value = 0.08011049723756906
$ws.Range("B3").Value = 0.004926108374384237
$ws.Range("C3").Value = 0.04433497536945813
$ws.Range("J3").Value = 0.01970443349753695
$ws.Range("P3").Value = 0.7536945812807881
$ws.Range("S3").Value = 0.1773399014778325
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.5813953488372093
$ws.Range("S4").Value = 0.3953488372093023
$ws.Range("B6").Value = 0.06481481481481481
$ws.Range("D6").Value = 0.01388888888888889
$ws.Range("F6").Value = 0.06481481481481481
$ws.Range("J6").Value = 0.2407407407407407
$ws.Range("O6").Value = 0.03240740740740741
$ws.Range("Q6").Value = 0.125
$ws.Range("R6").Value = 0.05092592592592592
$ws.Range("S6").Value = 0.4074074074074074
$ws.Range("B7").Value = 0.1273584905660377
$ws.Range("D7").Value = 0.01415094339622642
$ws.Range("F7").Value = 0.08018867924528301
$ws.Range("J7").Value = 0.1132075471698113
$ws.Range("O7").Value = 0.01415094339622642
$ws.Range("Q7").Value = 0.1509433962264151
$ws.Range("R7").Value = 0.09905660377358491
$ws.Range("S7").Value = 0.4009433962264151
$ws.Range("B8").Value = 0.1236203090507726
$ws.Range("D8").Value = 0.02207505518763797
$ws.Range("F8").Value = 0.06181015452538632
$ws.Range("J8").Value = 0.1169977924944812
$ws.Range("O8").Value = 0.02649006622516556
$ws.Range("Q8").Value = 0.1479028697571744
$ws.Range("R8").Value = 0.09271523178807947
$ws.Range("S8").Value = 0.4083885209713024
$ws.Range("B9").Value = 0.1071428571428571
$ws.Range("D9").Value = 0.02380952380952381
$ws.Range("F9").Value = 0.07738095238095238
$ws.Range("J9").Value = 0.1488095238095238
$ws.Range("O9").Value = 0.01785714285714286
$ws.Range("Q9").Value = 0.1130952380952381
$ws.Range("R9").Value = 0.05357142857142857
$ws.Range("S9").Value = 0.4583333333333333
$ws.Range("B10").Value = 0.1315996847911741
$ws.Range("D10").Value = 0.01891252955082742
$ws.Range("E10").Value = 0.0007880220646178094
$ws.Range("F10").Value = 0.06540583136327817
$ws.Range("J10").Value = 0.1252955082742317
$ws.Range("O10").Value = 0.02206461780929866
$ws.Range("Q10").Value = 0.2159180457052798
$ws.Range("R10").Value = 0.06461780929866036
$ws.Range("S10").Value = 0.355397951142632
$ws.Range("G11").Value = 0.1476510067114094
$ws.Range("J11").Value = 0.06711409395973154
$ws.Range("K11").Value = 0.2046979865771812
$ws.Range("L11").Value = 0.5570469798657718
$ws.Range("S11").Value = 0.02348993288590604
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1534090909090909
$ws.Range("K12").Value = 0.005681818181818182
$ws.Range("L12").Value = 0.0625
$ws.Range("S12").Value = 0.02840909090909091
$ws.Range("G13").Value = 0.8260869565217391
$ws.Range("J13").Value = 0.1521739130434783
$ws.Range("S13").Value = 0.02173913043478261
$ws.Range("F15").Value = 0.004878048780487805
$ws.Range("H15").Value = 0.1560975609756098
$ws.Range("I15").Value = 0.03902439024390244
$ws.Range("J15").Value = 0.3268292682926829
$ws.Range("K15").Value = 0.05853658536585366
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("N15").Value = 0.004878048780487805
$ws.Range("O15").Value = 0.06829268292682927
$ws.Range("S15").Value = 0.3170731707317073
$ws.Range("F16").Value = 0.01754385964912281
$ws.Range("H16").Value = 0.162280701754386
$ws.Range("I16").Value = 0.07894736842105263
$ws.Range("J16").Value = 0.4429824561403509
$ws.Range("K16").Value = 0.1184210526315789
$ws.Range("M16").Value = 0.0131578947368421
$ws.Range("N16").Value = 0.004385964912280702
$ws.Range("O16").Value = 0.04824561403508772
$ws.Range("S16").Value = 0.1140350877192982
$ws.Range("F17").Value = 0.01431980906921241
$ws.Range("H17").Value = 0.1646778042959427
$ws.Range("I17").Value = 0.1121718377088305
$ws.Range("J17").Value = 0.4391408114558473
$ws.Range("K17").Value = 0.09069212410501193
$ws.Range("M17").Value = 0.02147971360381861
$ws.Range("N17").Value = 0.002386634844868735
$ws.Range("O17").Value = 0.04534606205250596
$ws.Range("S17").Value = 0.1097852028639618
$ws.Range("F18").Value = 0.03105590062111801
$ws.Range("H18").Value = 0.1677018633540373
$ws.Range("I18").Value = 0.06832298136645963
$ws.Range("J18").Value = 0.4658385093167702
$ws.Range("K18").Value = 0.09937888198757763
$ws.Range("M18").Value = 0.006211180124223602
$ws.Range("O18").Value = 0.05590062111801242
$ws.Range("S18").Value = 0.1055900621118012
$ws.Range("F19").Value = 0.01494885916601102
$ws.Range("H19").Value = 0.2258064516129032
$ws.Range("I19").Value = 0.06845003933910307
$ws.Range("J19").Value = 0.3776553894571204
$ws.Range("K19").Value = 0.1132966168371361
$ws.Range("M19").Value = 0.02281667977970102
$ws.Range("N19").Value = 0.0007867820613690008
$ws.Range("O19").Value = 0.05900865460267506
$ws.Range("S19").Value = 0.1172305271439811
